$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$lo = $ws.ListObjects.Item(1)

# Current (before) layout of the "Variables" table:
#   D = variable-type
#   E = en_variable-label
#   F = en_note
#   G = en_elimination
#
# Target (after) layout:
#   D = variable-type        (unchanged)
#   E = elimination          (moved here from G, renamed from en_elimination)
#   F = en_variable-label    (moved here from E)
#   G = en_note              (moved here from F)
#
# The "en_elimination" flag column shouldn't be tied to a language, so it is
# renamed to "elimination" and relocated right after "variable-type".

Write-Output ("Table columns before: " + $lo.ListColumns.Item(5).Name + ", " + $lo.ListColumns.Item(6).Name + ", " + $lo.ListColumns.Item(7).Name)

# Capture the data that needs to move before it gets overwritten.
$varLabelHeader = $ws.Range("E1").Value2
$varLabelRow2 = $ws.Range("E2").Value2
$varLabelRow3 = $ws.Range("E3").Value2
$varLabelRow4 = $ws.Range("E4").Value2

$noteHeader = $ws.Range("F1").Value2
$noteRow2 = $ws.Range("F2").Value2
$noteRow3 = $ws.Range("F3").Value2
$noteRow4 = $ws.Range("F4").Value2

# Move "en_note" column to G (it currently has no data below the header).
$ws.Range("G1").Value = $noteHeader
$ws.Range("G2").Value = $noteRow2
$ws.Range("G3").Value = $noteRow3
$ws.Range("G4").Value = $noteRow4

# Move "en_variable-label" column to F.
$ws.Range("F1").Value = $varLabelHeader
$ws.Range("F2").Value = $varLabelRow2
$ws.Range("F3").Value = $varLabelRow3
$ws.Range("F4").Value = $varLabelRow4

# Put the renamed, language-independent "elimination" column in E (empty data body).
$ws.Range("E1").Value = "elimination"
$ws.Range("E2").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""

# The column widths travel with the data: what used to be column G's width is now
# column E's width, what used to be column E's width is now column F's, etc.
$ws.Range("E1").ColumnWidth = 17.65
$ws.Range("F1").ColumnWidth = 18.166666666666668
$ws.Range("G1").ColumnWidth = 9.83

# Update the view selection to match the new active cell.
$ws.Range("E2").Select()

Write-Output ("Table columns after: " + $lo.ListColumns.Item(5).Name + ", " + $lo.ListColumns.Item(6).Name + ", " + $lo.ListColumns.Item(7).Name)
